# Build site at 2021-01-29 15:15:03 UTC
# Apply the LOQ4240.xlsx content update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: Ativação date change ---
# Use a formula + paste-values round trip so the literal text "01/01/2021"
# is stored as a shared string instead of being auto-converted to a date serial.
$ws.Cells.Item(8, 2).Formula = '="01/01/2021"'
$ws.Cells.Item(8, 2).Copy() | Out-Null
$ws.Cells.Item(8, 2).PasteSpecial(-4163) | Out-Null
$ws.Cells.Item(8, 3).Formula = '="01/01/2021"'
$ws.Cells.Item(8, 3).Copy() | Out-Null
$ws.Cells.Item(8, 3).PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = $false

# --- Row 11: add English objectives text (new cells B11/C11) ---
# Copy wrap-text format from an existing formatted row before writing values.
$ws.Range("B13").Copy() | Out-Null
$ws.Range("B11").PasteSpecial(-4122) | Out-Null
$ws.Range("C13").Copy() | Out-Null
$ws.Range("C11").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(11, 2).Value = "Introduce the fundamental concepts of management science and organization selttings."
$ws.Cells.Item(11, 3).Value = "Introduce the fundamental concepts of management science and organization selttings."

# --- Row 13: Docente responsável change ---
$ws.Cells.Item(13, 2).Value = "11079086 - Herlandí de Souza Andrade"
$ws.Cells.Item(13, 3).Value = "11079086 - Herlandí de Souza Andrade"

# --- Row 14: Programa resumido text change (newline removed) ---
$ws.Cells.Item(14, 2).Value = "1. Áreas de Atuação da Administração.2. Estrutura organizacional."
$ws.Cells.Item(14, 3).Value = "1. Áreas de Atuação da Administração.2. Estrutura organizacional."

# --- Row 15: add English short syllabus text (new cells B15/C15) ---
$ws.Range("B14").Copy() | Out-Null
$ws.Range("B15").PasteSpecial(-4122) | Out-Null
$ws.Range("C14").Copy() | Out-Null
$ws.Range("C15").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(15, 2).Value = "1. Management Practice Areas. 2. Organizational structure"
$ws.Cells.Item(15, 3).Value = "1. Management Practice Areas. 2. Organizational structure"

# --- Row 16: Programa text change (newline removed) ---
$ws.Cells.Item(16, 2).Value = "1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização."
$ws.Cells.Item(16, 3).Value = "1. Noções básicas de Marketing, Finanças e Recursos Humanos. 2. Diferentes configurações de organização."

# --- Row 17: add English syllabus text (new cells B17/C17) ---
$ws.Range("B16").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Copy() | Out-Null
$ws.Range("C17").PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(17, 2).Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."
$ws.Cells.Item(17, 3).Value = "1. Basic notions of Marketing, Finance and Human Resources.2. Different organization settings."

# --- Row 19: Método text change ---
$ws.Cells.Item(19, 2).Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."
$ws.Cells.Item(19, 3).Value = "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

# --- Row 20: Critério text change ---
$ws.Cells.Item(20, 2).Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"
$ws.Cells.Item(20, 3).Value = "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"

# --- Row 21: Norma de recuperação text change ---
$ws.Cells.Item(21, 2).Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"
$ws.Cells.Item(21, 3).Value = "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"

# --- Row 22: Bibliografia text change ---
$ws.Cells.Item(22, 2).Value = "Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996."
$ws.Cells.Item(22, 3).Value = "Chiavenato, I. Gestão de Pessoas. 4 ed. São Paulo: Manole, 2014.Chiavenato, I. Recursos Humanos: o capital humano das organizações. 10 ed. Rio de Janeiro, Campus, 2015.ROBBINS, S. P.; DECENZO, D. A.; WOLTER, R. Fundamentos de Gestão de Pessoas. São Paulo, saraiva, 2013.KOTLER, P. - ARMSTRONG, G. Princípios De Marketing. 15 ed. São Paulo: Pearson, 2014.KOTLER, P.; KELLER, K. L. Administração De Marketing. 15 ed. São Paulo: Pearson, 2019.CHIAVENATO, I. Introdução À Teoria Geral da Administração. 9 ed. São Paulo: Manole, 2014. MAXIMIANO, A. C. A. Teoria Geral da Administração: da Revolução Urbana À Revolução Digital. 8 ed. São Paulo: Atlas, 2017.GUERRINI, F. M.; ESCRIÇÃO FILHO, E.; ROSIM, D. Administração Para Engenheiros. Rio de Janeiro: Campus, 2016.CHIAVENATO, I. Administração Para Não Administradores: a Gestão de Negócios Ao Alcance de Todos. 2 ed. São Paulo: Manole, 2011.SILVA, M. M. L. Administração para Estudantes e Profissionais de Áreas Técnicas. São Paulo: Brasport, 2018.GITMAN, L. J. - ZUTTER, C. J. Princípios de Administração Financeira. 14 ed. São Paulo: Perason, 2017.GROPPELLI, A. A.; NIKBAKHT, E. Administração Financeira. 3 ed. São Paulo: Saraiva, 2010.MARCOUSÉ, I.; SURRIDGE, M.; GILLESPIE, A. Finanças. São Paulo: Saraiva, 2013.BOLMAN, L.G.; DEAL, T.E. Reframing organizations. San Francisco, John Wiley, 2013KOTLER, P.. O Marketing sem segredos. 1 ed. Porto Alegre. Bookman, 2005MINTZBERG, H. Criando organizações eficazes. 2 ed. São Paulo, Atlas, 2006.MORGAN, G. Imagens da organização. São Paulo, Atlas, 1996."

# --- Rows 23-24: remove Requisitos section entirely ---
$ws.Rows("23:24").Delete() | Out-Null
